$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Forandrad) for all data rows 2-161: 46059 -> 46060
$ws.Range("C2:C161").Value = 46060

# Reorder rows 33-161 data (columns A, B, F, G) per target permutation
$ws.Range("A33").Value = "A 61986-2021"
$ws.Range("B33").Value = 44502
$ws.Range("F33").Value = "Kommuner"
$ws.Range("G33").Value = 1.7
$ws.Range("A34").Value = "A 42048-2021"
$ws.Range("B34").Value = 44426
$ws.Range("F34").Value = ""
$ws.Range("G34").Value = 0.7
$ws.Range("A35").Value = "A 62220-2022"
$ws.Range("B35").Value = 44922.66780092593
$ws.Range("F35").Value = ""
$ws.Range("G35").Value = 1.5
$ws.Range("A36").Value = "A 417-2023"
$ws.Range("B36").Value = 44929
$ws.Range("F36").Value = ""
$ws.Range("G36").Value = 1.7
$ws.Range("A37").Value = "A 22789-2022"
$ws.Range("B37").Value = 44715
$ws.Range("F37").Value = ""
$ws.Range("G37").Value = 6
$ws.Range("A38").Value = "A 29822-2023"
$ws.Range("B38").Value = 45107.56379629629
$ws.Range("F38").Value = ""
$ws.Range("G38").Value = 1.5
$ws.Range("A39").Value = "A 45054-2023"
$ws.Range("B39").Value = 45191
$ws.Range("F39").Value = ""
$ws.Range("G39").Value = 0.5
$ws.Range("A40").Value = "A 15684-2023"
$ws.Range("B40").Value = 45021
$ws.Range("F40").Value = ""
$ws.Range("G40").Value = 4.4
$ws.Range("A41").Value = "A 58150-2022"
$ws.Range("B41").Value = 44900.7941087963
$ws.Range("F41").Value = ""
$ws.Range("G41").Value = 2.5
$ws.Range("A42").Value = "A 10115-2023"
$ws.Range("B42").Value = 44986.34202546296
$ws.Range("F42").Value = "Sveaskog"
$ws.Range("G42").Value = 2.3
$ws.Range("A43").Value = "A 418-2023"
$ws.Range("B43").Value = 44929
$ws.Range("F43").Value = ""
$ws.Range("G43").Value = 0.9
$ws.Range("A44").Value = "A 419-2023"
$ws.Range("B44").Value = 44929
$ws.Range("F44").Value = ""
$ws.Range("G44").Value = 2.3
$ws.Range("A45").Value = "A 5564-2023"
$ws.Range("B45").Value = 44960
$ws.Range("F45").Value = ""
$ws.Range("G45").Value = 6.5
$ws.Range("A46").Value = "A 10235-2023"
$ws.Range("B46").Value = 44986.65542824074
$ws.Range("F46").Value = ""
$ws.Range("G46").Value = 2.9
$ws.Range("A47").Value = "A 420-2023"
$ws.Range("B47").Value = 44929
$ws.Range("F47").Value = ""
$ws.Range("G47").Value = 1.4
$ws.Range("A48").Value = "A 506-2023"
$ws.Range("B48").Value = 44930
$ws.Range("F48").Value = ""
$ws.Range("G48").Value = 1.7
$ws.Range("A49").Value = "A 48901-2023"
$ws.Range("B49").Value = 45209.52467592592
$ws.Range("F49").Value = ""
$ws.Range("G49").Value = 18.3
$ws.Range("A50").Value = "A 36868-2022"
$ws.Range("B50").Value = 44805.6289699074
$ws.Range("F50").Value = ""
$ws.Range("G50").Value = 1.6
$ws.Range("A51").Value = "A 25065-2025"
$ws.Range("B51").Value = 45799.70430555556
$ws.Range("F51").Value = ""
$ws.Range("G51").Value = 2.2
$ws.Range("A52").Value = "A 25063-2025"
$ws.Range("B52").Value = 45799.70287037037
$ws.Range("F52").Value = ""
$ws.Range("G52").Value = 1.5
$ws.Range("A53").Value = "A 60264-2024"
$ws.Range("B53").Value = 45642.8599537037
$ws.Range("F53").Value = ""
$ws.Range("G53").Value = 4.9
$ws.Range("A54").Value = "A 62221-2022"
$ws.Range("B54").Value = 44922
$ws.Range("F54").Value = ""
$ws.Range("G54").Value = 3.1
$ws.Range("A55").Value = "A 40687-2025"
$ws.Range("B55").Value = 45896.87452546296
$ws.Range("F55").Value = ""
$ws.Range("G55").Value = 3.3
$ws.Range("A56").Value = "A 40758-2025"
$ws.Range("B56").Value = 45897
$ws.Range("F56").Value = ""
$ws.Range("G56").Value = 5.1
$ws.Range("A57").Value = "A 40685-2025"
$ws.Range("B57").Value = 45896.87293981481
$ws.Range("F57").Value = ""
$ws.Range("G57").Value = 3.5
$ws.Range("A58").Value = "A 40686-2025"
$ws.Range("B58").Value = 45896.87371527778
$ws.Range("F58").Value = ""
$ws.Range("G58").Value = 2
$ws.Range("A59").Value = "A 41130-2025"
$ws.Range("B59").Value = 45898
$ws.Range("F59").Value = ""
$ws.Range("G59").Value = 2.3
$ws.Range("A60").Value = "A 28703-2025"
$ws.Range("B60").Value = 45820.23819444444
$ws.Range("F60").Value = ""
$ws.Range("G60").Value = 7.8
$ws.Range("A61").Value = "A 14383-2024"
$ws.Range("B61").Value = 45394
$ws.Range("F61").Value = ""
$ws.Range("G61").Value = 1.2
$ws.Range("A62").Value = "A 8877-2025"
$ws.Range("B62").Value = 45713.37665509259
$ws.Range("F62").Value = ""
$ws.Range("G62").Value = 7.3
$ws.Range("A63").Value = "A 53160-2023"
$ws.Range("B63").Value = 45229.37484953704
$ws.Range("F63").Value = ""
$ws.Range("G63").Value = 1.6
$ws.Range("A64").Value = "A 43124-2025"
$ws.Range("B64").Value = 45909.70793981481
$ws.Range("F64").Value = ""
$ws.Range("G64").Value = 2.3
$ws.Range("A65").Value = "A 43125-2025"
$ws.Range("B65").Value = 45909.72638888889
$ws.Range("F65").Value = ""
$ws.Range("G65").Value = 0.7
$ws.Range("A66").Value = "A 43013-2025"
$ws.Range("B66").Value = 45909.49075231481
$ws.Range("F66").Value = ""
$ws.Range("G66").Value = 1
$ws.Range("A67").Value = "A 29819-2025"
$ws.Range("B67").Value = 45825.84666666666
$ws.Range("F67").Value = ""
$ws.Range("G67").Value = 2.1
$ws.Range("A68").Value = "A 6355-2024"
$ws.Range("B68").Value = 45338.47413194444
$ws.Range("F68").Value = ""
$ws.Range("G68").Value = 1.9
$ws.Range("A69").Value = "A 35055-2024"
$ws.Range("B69").Value = 45527
$ws.Range("F69").Value = ""
$ws.Range("G69").Value = 0.7
$ws.Range("A70").Value = "A 58365-2022"
$ws.Range("B70").Value = 44901.62020833333
$ws.Range("F70").Value = ""
$ws.Range("G70").Value = 4.9
$ws.Range("A71").Value = "A 27671-2025"
$ws.Range("B71").Value = 45813
$ws.Range("F71").Value = ""
$ws.Range("G71").Value = 5.7
$ws.Range("A72").Value = "A 45248-2025"
$ws.Range("B72").Value = 45919.61190972223
$ws.Range("F72").Value = ""
$ws.Range("G72").Value = 1.4
$ws.Range("A73").Value = "A 45821-2025"
$ws.Range("B73").Value = 45923.59707175926
$ws.Range("F73").Value = ""
$ws.Range("G73").Value = 1.4
$ws.Range("A74").Value = "A 62219-2022"
$ws.Range("B74").Value = 44922.66564814815
$ws.Range("F74").Value = ""
$ws.Range("G74").Value = 1.4
$ws.Range("A75").Value = "A 61975-2021"
$ws.Range("B75").Value = 44502
$ws.Range("F75").Value = "Kommuner"
$ws.Range("G75").Value = 4.2
$ws.Range("A76").Value = "A 47649-2023"
$ws.Range("B76").Value = 45203
$ws.Range("F76").Value = ""
$ws.Range("G76").Value = 0.6
$ws.Range("A77").Value = "A 2761-2023"
$ws.Range("B77").Value = 44944.61564814814
$ws.Range("F77").Value = ""
$ws.Range("G77").Value = 3.6
$ws.Range("A78").Value = "A 31677-2025"
$ws.Range("B78").Value = 45834.37809027778
$ws.Range("F78").Value = ""
$ws.Range("G78").Value = 1
$ws.Range("A79").Value = "A 48279-2025"
$ws.Range("B79").Value = 45933.58763888889
$ws.Range("F79").Value = ""
$ws.Range("G79").Value = 2.1
$ws.Range("A80").Value = "A 47903-2025"
$ws.Range("B80").Value = 45932
$ws.Range("F80").Value = ""
$ws.Range("G80").Value = 1.6
$ws.Range("A81").Value = "A 19306-2025"
$ws.Range("B81").Value = 45769.56755787037
$ws.Range("F81").Value = "Sveaskog"
$ws.Range("G81").Value = 1.2
$ws.Range("A82").Value = "A 8969-2023"
$ws.Range("B82").Value = 44979.50427083333
$ws.Range("F82").Value = ""
$ws.Range("G82").Value = 0.6
$ws.Range("A83").Value = "A 15409-2024"
$ws.Range("B83").Value = 45401.34701388889
$ws.Range("F83").Value = ""
$ws.Range("G83").Value = 3
$ws.Range("A84").Value = "A 64438-2023"
$ws.Range("B84").Value = 45280.68280092593
$ws.Range("F84").Value = ""
$ws.Range("G84").Value = 5.8
$ws.Range("A85").Value = "A 34419-2025"
$ws.Range("B85").Value = 45846.66137731481
$ws.Range("F85").Value = ""
$ws.Range("G85").Value = 1
$ws.Range("A86").Value = "A 34402-2025"
$ws.Range("B86").Value = 45846.61403935185
$ws.Range("F86").Value = ""
$ws.Range("G86").Value = 5.3
$ws.Range("A87").Value = "A 34411-2025"
$ws.Range("B87").Value = 45846.64509259259
$ws.Range("F87").Value = ""
$ws.Range("G87").Value = 0.6
$ws.Range("A88").Value = "A 41931-2024"
$ws.Range("B88").Value = 45561.56418981482
$ws.Range("F88").Value = ""
$ws.Range("G88").Value = 5.1
$ws.Range("A89").Value = "A 51117-2025"
$ws.Range("B89").Value = 45947.52902777777
$ws.Range("F89").Value = ""
$ws.Range("G89").Value = 3.7
$ws.Range("A90").Value = "A 34530-2024"
$ws.Range("B90").Value = 45525
$ws.Range("F90").Value = ""
$ws.Range("G90").Value = 1
$ws.Range("A91").Value = "A 45665-2021"
$ws.Range("B91").Value = 44441.32771990741
$ws.Range("F91").Value = ""
$ws.Range("G91").Value = 1.3
$ws.Range("A92").Value = "A 10836-2023"
$ws.Range("B92").Value = 44991
$ws.Range("F92").Value = ""
$ws.Range("G92").Value = 2.6
$ws.Range("A93").Value = "A 10737-2023"
$ws.Range("B93").Value = 44984
$ws.Range("F93").Value = ""
$ws.Range("G93").Value = 3.4
$ws.Range("A94").Value = "A 52631-2024"
$ws.Range("B94").Value = 45610.30113425926
$ws.Range("F94").Value = ""
$ws.Range("G94").Value = 1.7
$ws.Range("A95").Value = "A 29648-2025"
$ws.Range("B95").Value = 45825
$ws.Range("F95").Value = ""
$ws.Range("G95").Value = 1.9
$ws.Range("A96").Value = "A 53651-2025"
$ws.Range("B96").Value = 45960.52734953703
$ws.Range("F96").Value = ""
$ws.Range("G96").Value = 4.4
$ws.Range("A97").Value = "A 54449-2025"
$ws.Range("B97").Value = 45965.54233796296
$ws.Range("F97").Value = ""
$ws.Range("G97").Value = 2.6
$ws.Range("A98").Value = "A 58101-2023"
$ws.Range("B98").Value = 45249.38699074074
$ws.Range("F98").Value = ""
$ws.Range("G98").Value = 6.1
$ws.Range("A99").Value = "A 16921-2025"
$ws.Range("B99").Value = 45755.37600694445
$ws.Range("F99").Value = "Sveaskog"
$ws.Range("G99").Value = 1.6
$ws.Range("A100").Value = "A 16618-2025"
$ws.Range("B100").Value = 45754.31761574074
$ws.Range("F100").Value = ""
$ws.Range("G100").Value = 0.8
$ws.Range("A101").Value = "A 56101-2025"
$ws.Range("B101").Value = 45973.71033564815
$ws.Range("F101").Value = ""
$ws.Range("G101").Value = 5
$ws.Range("A102").Value = "A 60254-2024"
$ws.Range("B102").Value = 45642
$ws.Range("F102").Value = ""
$ws.Range("G102").Value = 13.4
$ws.Range("A103").Value = "A 60260-2024"
$ws.Range("B103").Value = 45642
$ws.Range("F103").Value = ""
$ws.Range("G103").Value = 4.7
$ws.Range("A104").Value = "A 57527-2025"
$ws.Range("B104").Value = 45980.63293981482
$ws.Range("F104").Value = ""
$ws.Range("G104").Value = 3.2
$ws.Range("A105").Value = "A 53253-2023"
$ws.Range("B105").Value = 45229.52552083333
$ws.Range("F105").Value = ""
$ws.Range("G105").Value = 0.5
$ws.Range("A106").Value = "A 39459-2023"
$ws.Range("B106").Value = 45166.81715277778
$ws.Range("F106").Value = ""
$ws.Range("G106").Value = 1.7
$ws.Range("A107").Value = "A 55788-2023"
$ws.Range("B107").Value = 45239
$ws.Range("F107").Value = ""
$ws.Range("G107").Value = 2.2
$ws.Range("A108").Value = "A 19310-2025"
$ws.Range("B108").Value = 45769.57109953704
$ws.Range("F108").Value = "Sveaskog"
$ws.Range("G108").Value = 1.9
$ws.Range("A109").Value = "A 2763-2023"
$ws.Range("B109").Value = 44944.62079861111
$ws.Range("F109").Value = ""
$ws.Range("G109").Value = 2.5
$ws.Range("A110").Value = "A 35048-2024"
$ws.Range("B110").Value = 45527.64537037037
$ws.Range("F110").Value = ""
$ws.Range("G110").Value = 2.8
$ws.Range("A111").Value = "A 4316-2024"
$ws.Range("B111").Value = 45324.65905092593
$ws.Range("F111").Value = ""
$ws.Range("G111").Value = 0.6
$ws.Range("A112").Value = "A 61133-2025"
$ws.Range("B112").Value = 46000.49719907407
$ws.Range("F112").Value = ""
$ws.Range("G112").Value = 1.1
$ws.Range("A113").Value = "A 16924-2025"
$ws.Range("B113").Value = 45755.3825462963
$ws.Range("F113").Value = "Sveaskog"
$ws.Range("G113").Value = 5.3
$ws.Range("A114").Value = "A 61701-2025"
$ws.Range("B114").Value = 46002.53532407407
$ws.Range("F114").Value = ""
$ws.Range("G114").Value = 4.8
$ws.Range("A115").Value = "A 4002-2026"
$ws.Range("B115").Value = 46044
$ws.Range("F115").Value = ""
$ws.Range("G115").Value = 0.7
$ws.Range("A116").Value = "A 7179-2023"
$ws.Range("B116").Value = 44965
$ws.Range("F116").Value = ""
$ws.Range("G116").Value = 1.7
$ws.Range("A117").Value = "A 40671-2024"
$ws.Range("B117").Value = 45558.36265046296
$ws.Range("F117").Value = ""
$ws.Range("G117").Value = 2.1
$ws.Range("A118").Value = "A 21738-2023"
$ws.Range("B118").Value = 45063
$ws.Range("F118").Value = ""
$ws.Range("G118").Value = 1.4
$ws.Range("A119").Value = "A 45669-2023"
$ws.Range("B119").Value = 45194.87333333334
$ws.Range("F119").Value = ""
$ws.Range("G119").Value = 1.2
$ws.Range("A120").Value = "A 33522-2024"
$ws.Range("B120").Value = 45519
$ws.Range("F120").Value = ""
$ws.Range("G120").Value = 3.2
$ws.Range("A121").Value = "A 20430-2025"
$ws.Range("B121").Value = 45775.47766203704
$ws.Range("F121").Value = ""
$ws.Range("G121").Value = 4.5
$ws.Range("A122").Value = "A 45588-2022"
$ws.Range("B122").Value = 44845.54098379629
$ws.Range("F122").Value = ""
$ws.Range("G122").Value = 1.9
$ws.Range("A123").Value = "A 7310-2026"
$ws.Range("B123").Value = 46058.66219907408
$ws.Range("F123").Value = ""
$ws.Range("G123").Value = 2.3
$ws.Range("A124").Value = "A 17150-2023"
$ws.Range("B124").Value = 45034
$ws.Range("F124").Value = ""
$ws.Range("G124").Value = 1.3
$ws.Range("A125").Value = "A 65114-2021"
$ws.Range("B125").Value = 44515.41678240741
$ws.Range("F125").Value = ""
$ws.Range("G125").Value = 1.5
$ws.Range("A126").Value = "A 62223-2022"
$ws.Range("B126").Value = 44922.67252314815
$ws.Range("F126").Value = ""
$ws.Range("G126").Value = 0.4
$ws.Range("A127").Value = "A 1878-2022"
$ws.Range("B127").Value = 44574
$ws.Range("F127").Value = ""
$ws.Range("G127").Value = 4
$ws.Range("A128").Value = "A 39211-2023"
$ws.Range("B128").Value = 45166.41819444444
$ws.Range("F128").Value = ""
$ws.Range("G128").Value = 1.6
$ws.Range("A129").Value = "A 38792-2024"
$ws.Range("B129").Value = 45547.48993055556
$ws.Range("F129").Value = ""
$ws.Range("G129").Value = 2
$ws.Range("A130").Value = "A 59273-2023"
$ws.Range("B130").Value = 45253.55516203704
$ws.Range("F130").Value = ""
$ws.Range("G130").Value = 0.7
$ws.Range("A131").Value = "A 45881-2024"
$ws.Range("B131").Value = 45580
$ws.Range("F131").Value = ""
$ws.Range("G131").Value = 1.9
$ws.Range("A132").Value = "A 17611-2025"
$ws.Range("B132").Value = 45757.79642361111
$ws.Range("F132").Value = ""
$ws.Range("G132").Value = 2.7
$ws.Range("A133").Value = "A 19384-2023"
$ws.Range("B133").Value = 45049
$ws.Range("F133").Value = ""
$ws.Range("G133").Value = 1.3
$ws.Range("A134").Value = "A 16687-2025"
$ws.Range("B134").Value = 45754.45372685185
$ws.Range("F134").Value = ""
$ws.Range("G134").Value = 1.9
$ws.Range("A135").Value = "A 60258-2024"
$ws.Range("B135").Value = 45642
$ws.Range("F135").Value = ""
$ws.Range("G135").Value = 1.3
$ws.Range("A136").Value = "A 60259-2024"
$ws.Range("B136").Value = 45642
$ws.Range("F136").Value = ""
$ws.Range("G136").Value = 0.9
$ws.Range("A137").Value = "A 47650-2023"
$ws.Range("B137").Value = 45203
$ws.Range("F137").Value = ""
$ws.Range("G137").Value = 0.3
$ws.Range("A138").Value = "A 6872-2022"
$ws.Range("B138").Value = 44602.75136574074
$ws.Range("F138").Value = ""
$ws.Range("G138").Value = 0.8
$ws.Range("A139").Value = "A 58633-2023"
$ws.Range("B139").Value = 45251.56759259259
$ws.Range("F139").Value = ""
$ws.Range("G139").Value = 0.7
$ws.Range("A140").Value = "A 29728-2022"
$ws.Range("B140").Value = 44755.45465277778
$ws.Range("F140").Value = ""
$ws.Range("G140").Value = 2
$ws.Range("A141").Value = "A 17335-2025"
$ws.Range("B141").Value = 45756
$ws.Range("F141").Value = ""
$ws.Range("G141").Value = 4.7
$ws.Range("A142").Value = "A 48510-2024"
$ws.Range("B142").Value = 45591
$ws.Range("F142").Value = ""
$ws.Range("G142").Value = 0.5
$ws.Range("A143").Value = "A 54948-2024"
$ws.Range("B143").Value = 45619
$ws.Range("F143").Value = ""
$ws.Range("G143").Value = 3.5
$ws.Range("A144").Value = "A 9954-2025"
$ws.Range("B144").Value = 45719.35216435185
$ws.Range("F144").Value = ""
$ws.Range("G144").Value = 1
$ws.Range("A145").Value = "A 4241-2023"
$ws.Range("B145").Value = 44953
$ws.Range("F145").Value = ""
$ws.Range("G145").Value = 2.9
$ws.Range("A146").Value = "A 4242-2023"
$ws.Range("B146").Value = 44953
$ws.Range("F146").Value = ""
$ws.Range("G146").Value = 2.6
$ws.Range("A147").Value = "A 60253-2024"
$ws.Range("B147").Value = 45642
$ws.Range("F147").Value = ""
$ws.Range("G147").Value = 2.3
$ws.Range("A148").Value = "A 53179-2023"
$ws.Range("B148").Value = 45229.40875
$ws.Range("F148").Value = ""
$ws.Range("G148").Value = 4.5
$ws.Range("A149").Value = "A 42016-2023"
$ws.Range("B149").Value = 45177.47246527778
$ws.Range("F149").Value = ""
$ws.Range("G149").Value = 2.1
$ws.Range("A150").Value = "A 22028-2025"
$ws.Range("B150").Value = 45785.21606481481
$ws.Range("F150").Value = ""
$ws.Range("G150").Value = 1.2
$ws.Range("A151").Value = "A 3572-2025"
$ws.Range("B151").Value = 45680.71428240741
$ws.Range("F151").Value = ""
$ws.Range("G151").Value = 1.5
$ws.Range("A152").Value = "A 54412-2024"
$ws.Range("B152").Value = 45617
$ws.Range("F152").Value = ""
$ws.Range("G152").Value = 1.9
$ws.Range("A153").Value = "A 22418-2025"
$ws.Range("B153").Value = 45786.56079861111
$ws.Range("F153").Value = "Kyrkan"
$ws.Range("G153").Value = 0.7
$ws.Range("A154").Value = "A 22358-2025"
$ws.Range("B154").Value = 45786.46289351852
$ws.Range("F154").Value = "Kyrkan"
$ws.Range("G154").Value = 1.7
$ws.Range("A155").Value = "A 1191-2025"
$ws.Range("B155").Value = 45666
$ws.Range("F155").Value = ""
$ws.Range("G155").Value = 2.1
$ws.Range("A156").Value = "A 22828-2025"
$ws.Range("B156").Value = 45789.84356481482
$ws.Range("F156").Value = ""
$ws.Range("G156").Value = 0.6
$ws.Range("A157").Value = "A 12954-2025"
$ws.Range("B157").Value = 45734.40412037037
$ws.Range("F157").Value = ""
$ws.Range("G157").Value = 1.3
$ws.Range("A158").Value = "A 22827-2025"
$ws.Range("B158").Value = 45789.84092592593
$ws.Range("F158").Value = ""
$ws.Range("G158").Value = 1
$ws.Range("A159").Value = "A 45599-2024"
$ws.Range("B159").Value = 45579.43138888889
$ws.Range("F159").Value = ""
$ws.Range("G159").Value = 6.8
$ws.Range("A160").Value = "A 8450-2025"
$ws.Range("B160").Value = 45709.45519675926
$ws.Range("F160").Value = ""
$ws.Range("G160").Value = 1.4
$ws.Range("A161").Value = "A 10462-2025"
$ws.Range("B161").Value = 45720.7705787037
$ws.Range("F161").Value = ""
$ws.Range("G161").Value = 2.4
